$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Hunk 1: the empty (bold, size-28) paragraph right after the title
# becomes a hyperlink ("https://www.freecodecamp.org/learn/") followed
# by a space, a Wingdings arrow glyph and " free code camp"; the
# paragraph-mark formatting (b/bCs/sz/szCs/lang) that used to live in
# the empty <w:pPr><w:rPr> is dropped once real content is typed in,
# and the document's "_GoBack" edit-tracking bookmark moves from the
# end of the document to the end of this newly edited paragraph.
# ---------------------------------------------------------------------

# Move the _GoBack bookmark: delete its old (end-of-document) location;
# it gets re-created (below) right where the new text is typed.
$goBack = $d.Bookmarks.Item("_GoBack")
$goBack.Delete()

$p2 = $d.Paragraphs(2)
$r2 = $p2.Range
$newPara1 = @"
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships"><w:hyperlink r:id="rIdPLACEHOLDER" w:history="1"><w:r><w:rPr><w:rStyle w:val="Hyperlink"/></w:rPr><w:t>https://www.freecodecamp.org/learn/</w:t></w:r></w:hyperlink><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:sym w:font="Wingdings" w:char="F0E0"/></w:r><w:r><w:t xml:space="preserve"> free code camp</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p>
"@
$r2.InsertXML($newPara1)

# Fix up the hyperlink relationship/address (InsertXML mints a fresh
# relationship id automatically; make sure the Address + Hyperlink
# character style are set correctly) - mirrors what Word does when you
# paste/insert a hyperlink field.
$hl = $d.Hyperlinks.Item(2)
$hl.Address = "https://www.freecodecamp.org/learn/"
$hl.Range.Style = "Hyperlink"

# ---------------------------------------------------------------------
# Hunk 2: the cached "last rendered page break" marker moves from the
# start of the "Things to do later:" run to the start of the
# "1. stop the game after user wins" run (two paragraphs later) - pure
# Word re-pagination bookkeeping that happens on save.
# ---------------------------------------------------------------------

$todoPara = $null
$stopPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $txt = $d.Paragraphs($i).Range.Text.TrimEnd("`r")
    if ($txt -eq "Things to do later:") {
        $next = $d.Paragraphs($i + 1).Range.Text
        if ($next -like "*stop the game after user wins*") {
            $todoPara = $i
            $stopPara = $i + 1
            break
        }
    }
}

$p45 = $d.Paragraphs($todoPara)
$xml45 = @"
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>Things to do later:</w:t></w:r></w:p>
"@
$p45.Range.InsertXML($xml45)

$p46 = $d.Paragraphs($stopPara)
$xml46 = @"
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:lastRenderedPageBreak/><w:tab/><w:t>1. stop the game after user wins</w:t></w:r><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>.</w:t></w:r></w:p>
"@
$p46.Range.InsertXML($xml46)

Write-Host "edits applied"
